# Update cryptos list values per upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.966.97'
$ws.Range('D3').Value = '2.418.46'
$ws.Range('E3').Value = '  -0.16%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '562.56'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.13%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.92'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.94%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  -0.23%  '
$ws.Range('E9').Value = '  -0.13%  '
$ws.Range('E10').Value = '  -0.82%  '
$ws.Range('E11').Value = '  -4.04%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.349'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.71%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '26.19'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.30%  '
$ws.Range('E14').Value = '  -1.82%  '
$ws.Range('D15').Value = '2.847.89'
$ws.Range('E15').Value = '  -0.42%  '
$ws.Range('D16').Value = '61.954.94'
$ws.Range('E16').Value = '  -0.10%  '
$ws.Range('D17').Value = '2.393.51'
$ws.Range('E17').Value = '  -1.14%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.33'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.25%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '323.51'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.35%  '
$ws.Range('E20').Value = '  -1.20%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.82'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.92%  '
$ws.Range('E22').Value = '  -0.10%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '66.72'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.01%  '
$ws.Range('E24').Value = '  +0.69%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.77'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.80%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '552.14'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -5.88%  '
$ws.Range('E27').Value = '  +0.30%  '
$ws.Range('D28').Value = '2.536.19'
$ws.Range('E28').Value = '  +0.34%  '
$ws.Range('D29').Value = '0.0₃0932'
$ws.Range('E29').Value = '  -1.11%  '
$ws.Range('E30').Value = '  -0.51%  '
$ws.Range('E31').Value = '  -4.06%  '
$ws.Range('E32').Value = '  -1.99%  '
$ws.Range('E33').Value = '  -0.18%  '
$ws.Range('E34').Value = '  -3.45%  '
$ws.Range('E35').Value = '  -0.05%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.72'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.93%  '
$ws.Range('E37').Value = '  -1.68%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '152.59'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.09%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.42'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.94%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.55'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.64%  '
$ws.Range('E41').Value = '  -0.71%  '
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '147.18'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.20%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.22'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.81%  '
$ws.Range('E45').Value = '  -0.29%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0527'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.08%  '
$ws.Range('E47').Value = '  +0.38%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '19.82'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.60%  '
$ws.Range('E49').Value = '  -0.43%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0227'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.59%  '
$ws.Range('B51').Value = 'WhiteBITCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '11.58'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.78%  '
